# Scheduled data refresh: update hard-coded market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) across all job
# sheets with freshly pulled values. These are plain data cells (no
# formulas anywhere in the workbook), so cells are simply overwritten
# with their new literal numbers.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 272.86667
$ws.Range("I28").Value = 220.92857
$ws.Range("K28").Value = 220.92857
$ws.Range("M28").Value = 264.07143
$ws.Range("H40").Value = 1694.6
$ws.Range("I40").Value = 1039.875
$ws.Range("J40").Value = 2442.8572
$ws.Range("K40").Value = 1039.875
$ws.Range("L40").Value = 2442.8572
$ws.Range("M40").Value = -864.875
$ws.Range("N40").Value = -2792.8572
$ws.Range("H62").Value = 2416.6667
$ws.Range("I62").Value = 2388.889
$ws.Range("J62").Value = 2458.3333
$ws.Range("K62").Value = 2388.889
$ws.Range("L62").Value = 2458.3333
$ws.Range("M62").Value = -1764.889
$ws.Range("N62").Value = -3706.3333
$ws.Range("H64").Value = 4424.9375
$ws.Range("I64").Value = 4199.75
$ws.Range("K64").Value = 4199.75
$ws.Range("M64").Value = -3951.75
$ws.Range("H65").Value = 2416.6667
$ws.Range("I65").Value = 2388.889
$ws.Range("J65").Value = 2458.3333
$ws.Range("K65").Value = 11944.445
$ws.Range("L65").Value = 12291.6665
$ws.Range("M65").Value = -8824.445
$ws.Range("N65").Value = -18531.6665
$ws.Range("H67").Value = 4424.9375
$ws.Range("I67").Value = 4199.75
$ws.Range("K67").Value = 4199.75
$ws.Range("M67").Value = -3341.75
$ws.Range("H132").Value = 11171.429
$ws.Range("I132").Value = 13440
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 40320
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -37790
$ws.Range("N132").Value = -21560
$ws.Range("H138").Value = 1804.3846
$ws.Range("I138").Value = 1181.0312
$ws.Range("J138").Value = 2408.8484
$ws.Range("K138").Value = 3543.0936
$ws.Range("L138").Value = 7226.5452
$ws.Range("M138").Value = 1596.9064
$ws.Range("N138").Value = -17506.5452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1732.3182
$ws.Range("I32").Value = 1732.3182
$ws.Range("K32").Value = 1732.3182
$ws.Range("M32").Value = -1445.3182
$ws.Range("H74").Value = 55558560
$ws.Range("I74").Value = 66670060
$ws.Range("J74").Value = 1066.3334
$ws.Range("K74").Value = 66670060
$ws.Range("L74").Value = 1066.3334
$ws.Range("M74").Value = -66669186
$ws.Range("N74").Value = -2814.3334
$ws.Range("H77").Value = 55558560
$ws.Range("I77").Value = 66670060
$ws.Range("J77").Value = 1066.3334
$ws.Range("K77").Value = 333350300
$ws.Range("L77").Value = 5331.666999999999
$ws.Range("M77").Value = -333345932
$ws.Range("N77").Value = -14067.667
$ws.Range("H122").Value = 2085.7896
$ws.Range("I122").Value = 2385.8462
$ws.Range("J122").Value = 1435.6666
$ws.Range("K122").Value = 7157.5386
$ws.Range("L122").Value = 4306.9998
$ws.Range("M122").Value = -4707.5386
$ws.Range("N122").Value = -9206.9998
$ws.Range("H132").Value = 13613.523
$ws.Range("I132").Value = 1540.7858
$ws.Range("K132").Value = 4622.357400000001
$ws.Range("M132").Value = -2092.357400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3800.0908
$ws.Range("J20").Value = 1749.5
$ws.Range("L20").Value = 1749.5
$ws.Range("N20").Value = -2243.5
$ws.Range("H94").Value = 1101.24
$ws.Range("I94").Value = 955.55
$ws.Range("K94").Value = 955.55
$ws.Range("M94").Value = -504.55
$ws.Range("H99").Value = 2113.75
$ws.Range("I99").Value = 1303.3334
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 1303.3334
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = 194.6666
$ws.Range("N99").Value = -5596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 28478.572
$ws.Range("J74").Value = 30725
$ws.Range("L74").Value = 30725
$ws.Range("N74").Value = -32473
$ws.Range("H77").Value = 28478.572
$ws.Range("J77").Value = 30725
$ws.Range("L77").Value = 92175
$ws.Range("N77").Value = -100911
$ws.Range("H132").Value = 2403.366
$ws.Range("I132").Value = 1803.2667
$ws.Range("J132").Value = 4040
$ws.Range("K132").Value = 5409.800099999999
$ws.Range("L132").Value = 12120
$ws.Range("M132").Value = -2879.800099999999
$ws.Range("N132").Value = -17180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1257
$ws.Range("I5").Value = 789.125
$ws.Range("K5").Value = 2367.375
$ws.Range("M5").Value = -2255.375
$ws.Range("H114").Value = 1022.25
$ws.Range("I114").Value = 1419.5714
$ws.Range("K114").Value = 4258.7142
$ws.Range("M114").Value = -1004.7142
$ws.Range("H122").Value = 760.4286
$ws.Range("J122").Value = 1048.25
$ws.Range("L122").Value = 9434.25
$ws.Range("N122").Value = -14334.25
$ws.Range("H131").Value = 728.5599999999999
$ws.Range("J131").Value = 728.5599999999999
$ws.Range("L131").Value = 2185.68
$ws.Range("N131").Value = -12265.68
$ws.Range("H135").Value = 1257
$ws.Range("I135").Value = 789.125
$ws.Range("K135").Value = 7102.125
$ws.Range("M135").Value = -4567.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 19998.5
$ws.Range("J103").Value = 19998.5
$ws.Range("L103").Value = 19998.5
$ws.Range("N103").Value = -22342.5
$ws.Range("H126").Value = 3775.6316
$ws.Range("I126").Value = 2862.96
$ws.Range("J126").Value = 5530.769
$ws.Range("K126").Value = 8588.880000000001
$ws.Range("L126").Value = 16592.307
$ws.Range("M126").Value = -6118.880000000001
$ws.Range("N126").Value = -21532.307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 16657.4
$ws.Range("J76").Value = 16657.4
$ws.Range("L76").Value = 16657.4
$ws.Range("N76").Value = -17333.4
$ws.Range("H79").Value = 16657.4
$ws.Range("J79").Value = 16657.4
$ws.Range("L79").Value = 16657.4
$ws.Range("N79").Value = -18997.4
$ws.Range("H100").Value = 2287.5
$ws.Range("I100").Value = 1993.625
$ws.Range("J100").Value = 2875.25
$ws.Range("K100").Value = 1993.625
$ws.Range("L100").Value = 2875.25
$ws.Range("M100").Value = -1452.625
$ws.Range("N100").Value = -3957.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 19259.4
$ws.Range("J82").Value = 19259.4
$ws.Range("L82").Value = 19259.4
$ws.Range("N82").Value = -20025.4
$ws.Range("H85").Value = 19259.4
$ws.Range("J85").Value = 19259.4
$ws.Range("L85").Value = 19259.4
$ws.Range("N85").Value = -21911.4
$ws.Range("H132").Value = 3197.6
$ws.Range("I132").Value = 1996.6666
$ws.Range("K132").Value = 5989.9998
$ws.Range("M132").Value = -3459.9998
$ws.Range("H133").Value = 50715
$ws.Range("J133").Value = 50715
$ws.Range("L133").Value = 50715
$ws.Range("N133").Value = -60835
